$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text in A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$texto = $ws1.Range("A1").Value()
$texto = $texto.Replace("1000 Bs = 7.58 = 31012.27 pesos", "1000 Bs = 7.52 = 30752.63 pesos")
$texto = $texto.Replace("31012.27 pesos = 7.56 = 961.42 Bs", "30752.63 pesos = 7.48 = 962.16 Bs")
$ws1.Range("A1").Value = $texto

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 133
$ws2.Range("O10").Value = 4090.1
$ws2.Range("N12").Value = 4108.9
$ws2.Range("O12").Value = 128.555
